$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => [new C value, new E value]
$updates = @{
    63  = @(14355, 36182831)
    81  = @(17433, 133950839)
    91  = @(151134, 482406764)
    92  = @(409105, 1595107955)
    93  = @(209568, 1308770817)
    94  = @(94193, 917466940)
    95  = @(50759, 932280796)
    96  = @(17268, 792060406)
    104 = @(135239, 272198155)
    111 = @(116, 7791001)
    115 = @(11695, 32962731)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("C$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}
